$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the RNA isolation date (column A, rnaDate) for all data rows to 10.03.20.
# A leading apostrophe forces Excel to store the value as literal text instead
# of auto-converting the mm.dd.yy-looking string into a date serial number.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = "'10.03.20"
}
# Drop the quote-prefix formatting picked up above so the cells end up with no
# explicit style, matching a plain text cell entered directly.
$ws.Range("A2:A27").ClearFormats()

# Update the view selection to match the saved state
$ws.Range("A3:A27").Select()
